{"js": "// The Pandoc \"skylighting\" syntax-highlighting character styles in\n// styles.xml had their <w:rPr> children in an order that doesn't match\n// the CT_RPr content model in wml.xsd (w:color before w:b/w:i). Fix the\n// ordering by re-touching the bold/italic flags on the affected styles\n// so the run-properties are re-emitted in schema order (b, i, ... color).\nconst boldStyles = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\n\nconst italicStyles = [\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n];\n\nconst targetNames = new Set([...boldStyles, ...italicStyles]);\n\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\n\nfor (const s of styles.items) {\n  s.load(\"nameLocal\");\n}\nawait context.sync();\n\nconst byName = {};\nfor (const s of styles.items) {\n  if (targetNames.has(s.nameLocal)) {\n    byName[s.nameLocal] = s;\n  }\n}\n\nfor (const name of boldStyles) {\n  byName[name].font.bold = true;\n}\nfor (const name of italicStyles) {\n  byName[name].font.italic = true;\n}\n\nawait context.sync();\n", "ps1": "# The Pandoc \"skylighting\" syntax-highlighting character styles in\n# styles.xml had their <w:rPr> children in an order that doesn't match\n# the CT_RPr content model in wml.xsd (w:color before w:b/w:i). Fix the\n# ordering by re-touching the bold/italic flags on the affected styles\n# so the run-properties are re-emitted in schema order (b, i, ... color).\n$d = $word.ActiveDocument\n\n$boldStyles = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\n$italicStyles = @(\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"InformationTok\",\n    \"WarningTok\"\n)\n\nforeach ($name in $boldStyles) {\n    $d.Styles.Item($name).Font.Bold = $true\n}\n\nforeach ($name in $italicStyles) {\n    $d.Styles.Item($name).Font.Italic = $true\n}\n"}
